# Generate Report for Handback
# Refresh the localization-status report: the handback files are now in
# sync with en-US, so update Status, the Latest Handback DateTime stamps,
# and clear the stale "version is not the latest" Error Detail message.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns (E2, F2) ---
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# --- zh-cn sheet ---
$zhcn.Range("C2").Value  = $newStatus
$zhcn.Range("K2").Value  = "2016-08-25 12:51:35"
$zhcn.Range("P2").Value  = ""

# --- de-de sheet ---
$dede.Range("C2").Value  = $newStatus
$dede.Range("K2").Value  = "2016-08-25 12:51:42"
$dede.Range("P2").Value  = ""

# --- Column widths: the report's columns auto-fit to the refreshed
#     (now shorter/longer) Status and Error Detail text. ---
$overview.Columns.Item(5).AutoFit()
$overview.Columns.Item(6).AutoFit()

$zhcn.Columns.Item(3).AutoFit()
$zhcn.Columns.Item(16).AutoFit()

$dede.Columns.Item(3).AutoFit()
$dede.Columns.Item(16).AutoFit()
